$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("s1")
Write-Host $ws.Name
